$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2000
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -1931
$ws.Range("N43").ClearContents()

$ws.Range("H51").Value = 3500

$ws.Range("H54").Value = 22500
$ws.Range("J54").Value = 22500
$ws.Range("L54").Value = 22500
$ws.Range("N54").Value = -23472

$ws.Range("H62").Value = 3414.7144
$ws.Range("I62").Value = 2226.25
$ws.Range("K62").Value = 2226.25
$ws.Range("M62").Value = -1602.25

$ws.Range("H65").Value = 3414.7144
$ws.Range("I65").Value = 2226.25
$ws.Range("K65").Value = 11131.25
$ws.Range("M65").Value = -8011.25

$ws.Range("H70").Value = 2800
$ws.Range("J70").Value = 3450
$ws.Range("L70").Value = 10350
$ws.Range("N70").Value = -10890

$ws.Range("H73").Value = 2800
$ws.Range("J73").Value = 3450
$ws.Range("L73").Value = 10350
$ws.Range("N73").Value = -12222

$ws.Range("H80").Value = 632.5
$ws.Range("I80").Value = 496.25
$ws.Range("J80").Value = 723.3333
$ws.Range("K80").Value = 1488.75
$ws.Range("L80").Value = 2169.9999
$ws.Range("M80").Value = -490.75
$ws.Range("N80").Value = -4165.9999

$ws.Range("H83").Value = 632.5
$ws.Range("I83").Value = 496.25
$ws.Range("J83").Value = 723.3333
$ws.Range("K83").Value = 4466.25
$ws.Range("L83").Value = 6509.9997
$ws.Range("M83").Value = 525.75
$ws.Range("N83").Value = -16493.9997

$ws.Range("H92").Value = 661.25
$ws.Range("I92").Value = 458.4
$ws.Range("K92").Value = 458.4
$ws.Range("M92").Value = 789.6

$ws.Range("H100").Value = 1020.2
$ws.Range("I100").Value = 1020.2
$ws.Range("K100").Value = 1020.2
$ws.Range("M100").Value = -479.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 53
$ws.Range("I5").Value = 48.3
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 48.3
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = 63.7
$ws.Range("N5").Value = -324

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

$ws.Range("H45").Value = 2962.1614
$ws.Range("I45").Value = 2708.5925
$ws.Range("K45").Value = 2708.5925
$ws.Range("M45").Value = -2331.5925

$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1788
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 1972.7
$ws.Range("I74").Value = 1990.2222
$ws.Range("K74").Value = 1990.2222
$ws.Range("M74").Value = -1116.2222

$ws.Range("H77").Value = 1972.7
$ws.Range("I77").Value = 1990.2222
$ws.Range("K77").Value = 9951.110999999999
$ws.Range("M77").Value = -5583.110999999999

$ws.Range("H102").Value = 811.4286
$ws.Range("I102").Value = 811.4286
$ws.Range("K102").Value = 811.4286
$ws.Range("M102").Value = 810.5714

$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3450
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 53
$ws.Range("I4").Value = 48.3
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 48.3
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = 66.7
$ws.Range("N4").Value = -330

$ws.Range("H99").Value = 2902.5
$ws.Range("I99").Value = 2902.5
$ws.Range("K99").Value = 2902.5
$ws.Range("M99").Value = -1404.5

$ws.Range("H107").Value = 7513.4287
$ws.Range("I107").Value = 7513.4287
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 7513.4287
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -5593.4287
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 24000000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 24000000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 24000000
$ws.Range("N6").Value = -24000226
$ws.Range("M6").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H31").Value = 15240.75
$ws.Range("I31").Value = 21595
$ws.Range("K31").Value = 21595
$ws.Range("M31").Value = -21300

$ws.Range("H34").Value = 15240.75
$ws.Range("I34").Value = 21595
$ws.Range("K34").Value = 21595
$ws.Range("M34").Value = -21393

$ws.Range("H41").Value = 10000
$ws.Range("I41").Value = 10000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 10000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -9572
$ws.Range("N41").ClearContents()

$ws.Range("H58").Value = 2955
$ws.Range("I58").Value = 2443.75
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 2443.75
$ws.Range("L58").Value = 5000
$ws.Range("M58").Value = -2240.75
$ws.Range("N58").Value = -5406

$ws.Range("H59").Value = 35500
$ws.Range("J59").Value = 35000
$ws.Range("L59").Value = 35000
$ws.Range("N59").Value = -37290

$ws.Range("H134").Value = 958
$ws.Range("I134").Value = 958
$ws.Range("K134").Value = 2874
$ws.Range("M134").Value = -339

$ws.Range("H136").Value = 2955
$ws.Range("I136").Value = 2443.75
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 7331.25
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -4781.25
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5087.25
$ws.Range("I80").Value = 1849
$ws.Range("J80").Value = 6166.6665
$ws.Range("K80").Value = 5547
$ws.Range("L80").Value = 18499.9995
$ws.Range("M80").Value = -4611
$ws.Range("N80").Value = -20371.9995

$ws.Range("H83").Value = 5087.25
$ws.Range("I83").Value = 1849
$ws.Range("J83").Value = 6166.6665
$ws.Range("K83").Value = 16641
$ws.Range("L83").Value = 55499.9985
$ws.Range("M83").Value = -11961
$ws.Range("N83").Value = -64859.9985

$ws.Range("H117").Value = 1777.6
$ws.Range("J117").Value = 2207.8333
$ws.Range("L117").Value = 6623.499899999999
$ws.Range("N117").Value = -13507.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 950.5
$ws.Range("I107").Value = 834
$ws.Range("K107").Value = 834
$ws.Range("M107").Value = 1086

$ws.Range("H113").Value = 955
$ws.Range("I113").Value = 893.75
$ws.Range("K113").Value = 893.75
$ws.Range("M113").Value = 1276.25

$ws.Range("H123").Value = 58663
$ws.Range("J123").Value = 58663
$ws.Range("L123").Value = 58663
$ws.Range("N123").Value = -63563

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H17").Value = 10000000
$ws.Range("I17").Value = 10000000
$ws.Range("K17").Value = 10000000
$ws.Range("M17").Value = -9999828

Write-Host "Edit complete"
